$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Template")

# --- Row 6: new electrofishing record, mirrors row 5's layout/style ---
# Copy cell-level formatting (number format / alignment / borders / fill / font)
# from the corresponding cell in row 5 so the new row reuses the existing
# style entries instead of minting new ones.
$fmtCells = @("A","B","C","D","E","F","L","M","N","P","R","S")
foreach ($col in $fmtCells) {
    $ws.Range($col + "5").Copy() | Out-Null
    $ws.Range($col + "6").PasteSpecial(-4122) | Out-Null
}
$excel.CutCopyMode = 0

# --- Values for row 6 ---
$ws.Range("A6").Value = 2021
$ws.Range("B6").Value = "Apr"
$ws.Range("C6").Value = 6
$ws.Range("E6").Value = "PWR"
$ws.Range("G6").Value = "WS"
$ws.Range("L6").Value = 45.549656159192402
$ws.Range("M6").Value = -65.013694691467194
$ws.Range("R6").Value = "QS"
$ws.Range("S6").Value = "QS"
$ws.Range("U6").Value = 150

# Row height matches the other data rows.
$ws.Rows.Item(6).RowHeight = 15.75

# Restore the active selection the workbook was left on.
$ws.Range("E17").Select() | Out-Null
